# "updated main GSC export data"
# Append two new daily data points to the "Chart" sheet (the GSC export's
# main data table): 2025-12-10 and 2025-12-11, each with Invalid=0 and
# Valid=29, continuing directly after the existing last row (2025-12-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

function Set-TextCell($cell, [string]$text) {
    # Column A holds plain text values that look like ISO dates
    # ("2025-12-10"). A bare string .Value assignment gets auto-coerced to a
    # real date serial by the "smart" input parser, which is not what the
    # source data uses (existing cells are plain shared-string text in
    # General format). Temporarily forcing Text format defeats the date
    # auto-detection; ClearFormats afterwards returns the cell to the
    # default (General) formatting used by every other row.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextCell $ws.Cells.Item($newRow, 1) "2025-12-10"
$ws.Cells.Item($newRow, 2).Value = 0.0
$ws.Cells.Item($newRow, 3).Value = 29.0

$newRow2 = $newRow + 1
Set-TextCell $ws.Cells.Item($newRow2, 1) "2025-12-11"
$ws.Cells.Item($newRow2, 2).Value = 0.0
$ws.Cells.Item($newRow2, 3).Value = 29.0
